# Add login and signup pages with form validation and styling
# (workbook-side bookkeeping that accompanied the app change)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("requirements")

# Requirement 2's description now reflects the Next.js API + Supabase backend
# (was: "Design and implement Express.js backend structure with RESTful routes")
$ws.Range("B3").Value = "Design and implement Next API + Supabase backend structure with RESTful routes"

# Requirement 2 moved from Sprint 1 to Sprint 2
$ws.Range("E3").Value = 2

# Sprint 2's requirement list now also includes requirement 2
$ws.Range("I3").Value = "2, 6, 7, 8, 9.1, 9.2"

# Update the active selection to match
$ws.Range("I3").Select()

# Remove the now-unused empty "sprints" and "Sheet1" worksheets
$wb.Worksheets.Item("sprints").Delete()
$wb.Worksheets.Item("Sheet1").Delete()
